# Append two new species-observation rows (17 and 18) to the "Artfynd" sheet,
# mirroring the rows already present (row 16 etc.) in shape/typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($ws, $row, $col, $value) {
    # Plain numeric cell.
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-TextCell($ws, $row, $col, $text) {
    # Force text storage even when the content looks like a number/date
    # (e.g. "1", "2023-09-07") by using Excel's leading-apostrophe
    # quote-prefix convention; the apostrophe itself is not stored.
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

function Set-BoolCell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Column layout (same columns populated as the existing data rows):
# A Id | B Taxonsorteringsordning | C Valideringsstatus | D Rödlistade
# E TaxonId | F Artnamn | G Vetenskapligt namn | H Auktor | I Antal
# J Enhet | K Ålder-Stadium | L Kön | M Aktivitet | N Metod | P Lokalnamn
# Q Ost | R Nord | S Noggrannhet | T Län | U Kommun | V Provins
# W Församling | Y Startdatum | Z Starttid | AA Slutdatum | AB Sluttid
# AD Ej återfunnen | AE Osäker artbestämning | AF Bestämningsmetod
# AG Ospontan | AT Bestämningsår | AW Rapportör | AX Observatörer
# AY Projektnamn

$rows = @(
    @{
        Row = 17
        A = 111957066; B = 55652; C = "Ovaliderad"; D = "LC"; E = 208255
        F = "Skogsödla"; G = "Zootoca vivipara"; H = "(Jacquin, 1787)"
        I = "1"; J = ""; K = "årsunge"; L = ""; M = ""; N = "observerad"
        P = "Sannahed sandtag Backar och Lund, Sannahed, Nrk"
        Q = 509365.3999776145; R = 6551082.436783144; S = 55
        T = "Örebro"; U = "Kumla"; V = "Närke"; W = "Kumla"
        Y = "2023-09-07"; Z = "00:00"; AA = "2023-09-07"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false; AT = ""
        AW = "David Bartholdsson"; AX = "David Bartholdsson"; AY = ""
    },
    @{
        Row = 18
        A = 111957069; B = 55652; C = "Ovaliderad"; D = "LC"; E = 208255
        F = "Skogsödla"; G = "Zootoca vivipara"; H = "(Jacquin, 1787)"
        I = "1"; J = ""; K = "adult"; L = ""; M = ""; N = "observerad"
        P = "Sannahed sandtag Äng, Sannahed, Nrk"
        Q = 509498.1630738945; R = 6551086.360099105; S = 75
        T = "Örebro"; U = "Kumla"; V = "Närke"; W = "Kumla"
        Y = "2023-09-07"; Z = "00:00"; AA = "2023-09-07"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false; AT = ""
        AW = "David Bartholdsson"; AX = "David Bartholdsson"; AY = ""
    }
)

# Column -> index, value-kind (num / text(forced) / bool)
$colSpec = @(
    @{ Col = "A";  Idx = 1;  Kind = "num" },
    @{ Col = "B";  Idx = 2;  Kind = "num" },
    @{ Col = "C";  Idx = 3;  Kind = "text" },
    @{ Col = "D";  Idx = 4;  Kind = "text" },
    @{ Col = "E";  Idx = 5;  Kind = "num" },
    @{ Col = "F";  Idx = 6;  Kind = "text" },
    @{ Col = "G";  Idx = 7;  Kind = "text" },
    @{ Col = "H";  Idx = 8;  Kind = "text" },
    @{ Col = "I";  Idx = 9;  Kind = "text" },
    @{ Col = "J";  Idx = 10; Kind = "text" },
    @{ Col = "K";  Idx = 11; Kind = "text" },
    @{ Col = "L";  Idx = 12; Kind = "text" },
    @{ Col = "M";  Idx = 13; Kind = "text" },
    @{ Col = "N";  Idx = 14; Kind = "text" },
    @{ Col = "P";  Idx = 16; Kind = "text" },
    @{ Col = "Q";  Idx = 17; Kind = "num" },
    @{ Col = "R";  Idx = 18; Kind = "num" },
    @{ Col = "S";  Idx = 19; Kind = "num" },
    @{ Col = "T";  Idx = 20; Kind = "text" },
    @{ Col = "U";  Idx = 21; Kind = "text" },
    @{ Col = "V";  Idx = 22; Kind = "text" },
    @{ Col = "W";  Idx = 23; Kind = "text" },
    @{ Col = "Y";  Idx = 25; Kind = "text" },
    @{ Col = "Z";  Idx = 26; Kind = "text" },
    @{ Col = "AA"; Idx = 27; Kind = "text" },
    @{ Col = "AB"; Idx = 28; Kind = "text" },
    @{ Col = "AD"; Idx = 30; Kind = "bool" },
    @{ Col = "AE"; Idx = 31; Kind = "bool" },
    @{ Col = "AF"; Idx = 32; Kind = "text" },
    @{ Col = "AG"; Idx = 33; Kind = "bool" },
    @{ Col = "AT"; Idx = 46; Kind = "text" },
    @{ Col = "AW"; Idx = 49; Kind = "text" },
    @{ Col = "AX"; Idx = 50; Kind = "text" },
    @{ Col = "AY"; Idx = 51; Kind = "text" }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row
    foreach ($spec in $colSpec) {
        $value = $rowData[$spec.Col]
        if ($spec.Kind -eq "num") {
            Set-NumCell $ws $r $spec.Idx $value
        } elseif ($spec.Kind -eq "bool") {
            Set-BoolCell $ws $r $spec.Idx $value
        } else {
            Set-TextCell $ws $r $spec.Idx $value
        }
    }
}
